$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - only B2 changes slightly
$ws.Range("B2").Value = 46640842766139.2

# Row 3: RandomForestRegressor - B3, C3, D3 change
$ws.Range("B3").Value = 33266025416571.03
$ws.Range("C3").Value = 38143163041203.87
$ws.Range("D3").Value = 36451670223412.15

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, with new values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 34344119739998.72
$ws.Range("C4").Value = 33242532460768.36
$ws.Range("D4").Value = 34344119739998.72

# Row 5: AdaBoostRegressor -> MLPRegressor, with new values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 50668038803691.88
$ws.Range("C5").Value = 39780181238528.34
$ws.Range("D5").Value = 33511962025915.6
